# Add the new "Balance After Transaction" column (E) to the Transactions
# header row, set the new column's width, update the dimension/selection,
# and leave the active selection on F2 (matching the authored template).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Balance After Transaction" column.
$ws.Range("E1").Value = "Balance After Transaction"

# Bold header style used by the other header cells (A1:D1) in this sheet.
$ws.Range("E1").Font.Bold = $true

# Match the authored column width for the new column (character units).
# Excel's ColumnWidth COM property and the stored XML "width" attribute
# differ by a constant 5/6 (~0.8333) offset on this workbook's font/DPI, so
# back that out to land on the authored width of 25.5 in the saved file.
$ws.Columns.Item(5).ColumnWidth = 24.666666666666668

# Move / collapse the selection to F2, as captured in the saved template.
$ws.Range("F2").Select()
